# Workbook was re-uploaded with: the A column values bumped from 3 to 4
# (rows 2-20), the sheet selection moved from C23 to A2:A20 (active cell
# A2), and the saved window height shrunk from 13660 to 11500.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bump A2:A20 from 3 -> 4 ---
$ws.Range("A2:A20").Value = 4

# --- Move the selection to A2:A20 (active cell A2) ---
$ws.Range("A2:A20").Select()

# --- Shrink the saved window height (32000 x 13660 -> 32000 x 11500) ---
$win = $excel.ActiveWindow
$win.Height = 11500
